# issue #5: stock data output to json file
#
# The "股票" (stock) worksheet gets a new "property_category" column
# inserted between the "total" and "date" columns, with the value
# "stock" for the existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before column H (which currently holds "date"),
# shifting date / legislator_name / legislator_id one column to the right.
$ws.Columns.Item(8).Insert()

# Populate the freshly inserted column with header + value.
$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
